$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. The sheet currently named "总计" (4th sheet) becomes "2022-Q1" and
#    is refilled with this quarter's fund-holding detail (same shape
#    as the 2021-Q4 sheet). A brand-new "总计" sheet is appended after
#    it, carrying the updated quarter-over-quarter summary table.
# ------------------------------------------------------------------

$qSheet = $wb.Worksheets.Item(4)
$qSheet.Name = "2022-Q1"

# Stretch the existing header/index styling (style already used by B1:D1
# and A2:A4) across the extra columns/row this sheet needs (E1:H1, A5)
# before writing their values, so the new cells pick up the same look.
$qSheet.Range("B1").Copy($qSheet.Range("E1:H1"))
$qSheet.Range("A4").Copy($qSheet.Range("A5"))

# -- Header row --
$qSheet.Range("B1").Value = "基金代码"
$qSheet.Range("C1").Value = "基金名称"
$qSheet.Range("D1").Value = "基金规模"
$qSheet.Range("E1").Value = "股票总仓位"
$qSheet.Range("F1").Value = "仓位占比"
$qSheet.Range("G1").Value = "持有市值(亿元)"
$qSheet.Range("H1").Value = "仓位排名"

# Columns B,C,D,E,F,G hold text that looks numeric (fund codes, money
# amounts) in the source data, so force text formatting before writing
# to avoid Excel silently coercing them (and dropping leading zeros).
$qSheet.Range("B2:G5").NumberFormat = "@"

# -- Row 2: 501201 --
$qSheet.Range("A2").Value = 0
$qSheet.Range("B2").Value = "501201"
$qSheet.Range("C2").Value = "红土创新科技创新 3 年封闭运作灵活配置混合"
$qSheet.Range("D2").Value = "3.99"
$qSheet.Range("E2").Value = "96.70"
$qSheet.Range("F2").Value = "3.28"
$qSheet.Range("G2").Value = "0.1309"
$qSheet.Range("H2").Value = 8

# -- Row 3: 010690 --
$qSheet.Range("A3").Value = 1
$qSheet.Range("B3").Value = "010690"
$qSheet.Range("C3").Value = "万家互联互通核心资产量化策略混合A"
$qSheet.Range("D3").Value = "0.85"
$qSheet.Range("E3").Value = "94.05"
$qSheet.Range("F3").Value = "5.15"
$qSheet.Range("G3").Value = "0.0438"
$qSheet.Range("H3").Value = 1

# -- Row 4: 168401 --
$qSheet.Range("A4").Value = 2
$qSheet.Range("B4").Value = "168401"
$qSheet.Range("C4").Value = "红土创新转型精选灵活配置混合（LOF）"
$qSheet.Range("D4").Value = "0.78"
$qSheet.Range("E4").Value = "93.82"
$qSheet.Range("F4").Value = "3.25"
$qSheet.Range("G4").Value = "0.0254"
$qSheet.Range("H4").Value = 9

# -- Row 5: 010691 --
$qSheet.Range("A5").Value = 3
$qSheet.Range("B5").Value = "010691"
$qSheet.Range("C5").Value = "万家互联互通核心资产量化策略混合C"
$qSheet.Range("D5").Value = "0.20"
$qSheet.Range("E5").Value = "94.05"
$qSheet.Range("F5").Value = "5.15"
$qSheet.Range("G5").Value = "0.0103"
$qSheet.Range("H5").Value = 1

# ------------------------------------------------------------------
# 2. Append the new "总计" summary sheet right after "2022-Q1", with
#    a fresh row for the new quarter on top of the previous totals.
# ------------------------------------------------------------------

$totalSheetTmp = $wb.Worksheets.Add()
$totalSheetTmp.Name = "总计"
$moveAnchor = $wb.Worksheets.Item("2022-Q1")
$totalSheetTmp.Move($null, $moveAnchor)
# Worksheet references go stale once the sheet collection reshuffles
# (rename/add/move all renumber indices), so re-resolve both sheets by
# name AFTER the move before using them any further.
$totalSheet = $wb.Worksheets.Item("总计")
$qSheet2 = $wb.Worksheets.Item("2022-Q1")

# Seed the header/index styling (style already used by B1:D1 and A2:A4
# on the "2022-Q1" sheet) onto this brand-new sheet before writing
# values. Copy from a SINGLE source cell each time — the COM host's
# Range.Copy only reliably fans a single source cell out across a
# multi-cell destination, not a multi-cell source range.
$qSheet2.Range("B1").Copy($totalSheet.Range("B1:D1"))
$qSheet2.Range("A2").Copy($totalSheet.Range("A2:A5"))

# -- Header row --
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# -- Row 2: 2022-Q1 (new) --
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.21

# -- Row 3: 2021-Q4 --
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 8
$totalSheet.Range("D3").Value = 2.34

# -- Row 4: 2021-Q3 --
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.16

# -- Row 5: 2021-Q2 --
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.11


